# Generate Report for Handoff
# Adds a new localization-status row for file 996f345d-76ed-4a36-aca8-7897bfd7bef4.md
# to the Overview sheet and to each language sheet (zh-cn, de-de), mirroring the
# existing 089084b7-a469-4f48-856b-98300d6a7fc4.md row that is already present.

$wb = $excel.ActiveWorkbook

$newGuidFile   = "996f345d-76ed-4a36-aca8-7897bfd7bef4.md"
$newGuid       = "996f345d-76ed-4a36-aca8-7897bfd7bef4"
$newXlfHash    = "94c73595a88e2263af8b912a54216edc1884bc38"

$readyStatus   = "Ready for handoff"
$includeText   = "Include"
$epochText     = "0001-01-01 00:00:00"

$mdCommit      = "85478c745bd0d5db2800773d6f23f9abd8345986"
$mdUrl         = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$newGuidFile"

# Hyperlink underline font color used by the workbook's "HyperLink" style (FF6495ED).
$hyperlinkColor = 15570276
$dateNumberFormat = "yyyy-mm-dd HH:mm:ss"

# Apply the hyperlink look (underline + blue) *after* Hyperlinks.Add, because
# Add() re-stamps its own built-in Hyperlink font on the range.
function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newGuidFile
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus
$wsOverview.Range("D3").Value = "2016-30-18 16:30:55"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, $null, $null, $newGuidFile) | Out-Null
Style-AsHyperlink $wsOverview.Range("A3")

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlfName = "$newGuid.$newXlfHash.zh-cn.xlf"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/29fb1cbc2de41897035ee53ffaf45721a48358e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"

$wsZhCn.Range("A3").Value = $newGuidFile
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $readyStatus
$wsZhCn.Range("D3").Value = $zhXlfName
$wsZhCn.Range("E3").NumberFormat = $dateNumberFormat
$wsZhCn.Range("E3").Value = "2016-03-18 16:30:52"
$wsZhCn.Range("H3").Value = $epochText
$wsZhCn.Range("I3").Value = $includeText

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdUrl, $null, $null, $newGuidFile) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), $mdUrl, $null, $null, ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $zhXlfUrl, $null, $null, $zhXlfName) | Out-Null
Style-AsHyperlink $wsZhCn.Range("A3")
Style-AsHyperlink $wsZhCn.Range("B3")
Style-AsHyperlink $wsZhCn.Range("D3")

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlfName = "$newGuid.$newXlfHash.de-de.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1eef84c32f57a4118fe30615b9e9781666cdfb24/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

$wsDeDe.Range("A3").Value = $newGuidFile
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $readyStatus
$wsDeDe.Range("D3").Value = $deXlfName
$wsDeDe.Range("E3").NumberFormat = $dateNumberFormat
$wsDeDe.Range("E3").Value = "2016-03-18 16:30:55"
$wsDeDe.Range("H3").Value = $epochText
$wsDeDe.Range("I3").Value = $includeText

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdUrl, $null, $null, $newGuidFile) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), $mdUrl, $null, $null, ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $deXlfUrl, $null, $null, $deXlfName) | Out-Null
Style-AsHyperlink $wsDeDe.Range("A3")
Style-AsHyperlink $wsDeDe.Range("B3")
Style-AsHyperlink $wsDeDe.Range("D3")

Write-Output "Added handoff rows for $newGuidFile to Overview, zh-cn, de-de"
